$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (C) column date for all existing data rows (2-27)
# from 45258 to 45259.
$ws.Range("C2:C27").Value = 45259

# Remove the two trailing rows (A59637-2023 / row 28 and A59686-2023 / row 29)
# that no longer belong in the sheet. Delete from bottom to top so row
# numbers of earlier rows are unaffected while deleting.
$ws.Rows.Item(29).Delete()
$ws.Rows.Item(28).Delete()

# Row 27 no longer carries an explicit custom row height in the target;
# AutoFit clears the stored height/customHeight attributes back to default.
$ws.Rows.Item(27).EntireRow.AutoFit()
